$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append (row, dateSerial, nuoviPos, sommaMobile7gg, sommaMobile7ggPer100k)
$data = @(
    @(465, 44539, 14, 46, 261.8248050543571),
    @(466, 44540, 12, 56, 318.7432409357391),
    @(467, 44541, 2, 54, 307.3595537594627),
    @(468, 44542, 15, 66, 375.6616768171211),
    @(469, 44543, 10, 68, 387.0453639933975),
    @(470, 44544, 10, 67, 381.3535204052593),
    @(471, 44545, 2, 65, 369.9698332289829),
    @(472, 44546, 8, 59, 335.8187717001537),
    @(473, 44547, 17, 64, 364.2779896408447),
    @(474, 44548, 13, 75, 426.8882691103648),
    @(475, 44550, 10, 70, 398.4290511696739),
    @(476, 44551, 7, 67, 381.3535204052593),
    @(477, 44552, 3, 60, 341.5106152882919),
    @(478, 44553, 12, 70, 398.4290511696739),
    @(479, 44554, 5, 67, 381.3535204052593),
    @(480, 44555, 13, 63, 358.5861460527065),
    @(481, 44556, 14, 64, 364.2779896408447),
    @(482, 44557, 35, 89, 506.5740793442996),
    @(483, 44558, 4, 86, 489.498548579885),
    @(484, 44559, 11, 94, 535.0332972849906),
    @(485, 44560, 26, 108, 614.7191075189254),
    @(486, 44561, 63, 166, 944.8460356309408),
    @(487, 44562, 42, 195, 1109.909499686949),
    @(488, 44563, 27, 208, 1183.903466332745),
    @(489, 44564, 27, 200, 1138.36871762764),
    @(490, 44565, 5, 201, 1144.060561215778),
    @(491, 44566, 10, 200, 1138.36871762764)
)

# Reference row that already has the right formatting (date style on col A)
$templateRow = 464

foreach ($entry in $data) {
    $r = $entry[0]
    $dateSerial = $entry[1]
    $nuoviPos = $entry[2]
    $sommaMobile = $entry[3]
    $sommaMobile100k = $entry[4]

    # Copy formatting (style) from the template row so new cells match existing ones
    $ws.Range("A$templateRow" + ":D$templateRow").Copy() | Out-Null
    $ws.Range("A$r" + ":D$r").PasteSpecial(-4122) | Out-Null # xlPasteFormats

    $ws.Cells.Item($r, 1).Value2 = $dateSerial
    $ws.Cells.Item($r, 2).Value2 = $nuoviPos
    $ws.Cells.Item($r, 3).Value2 = $sommaMobile
    $ws.Cells.Item($r, 4).Value2 = $sommaMobile100k
}

$excel.CutCopyMode = 0
